# Auto-generated edit script: applies the scheduled-runner market-price update
# to the per-item profit columns (H:N) across all eight class sheets.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 61.333332
$ws.Range("I9").Value = 61.333332
$ws.Range("K9").Value = 61.333332
$ws.Range("M9").Value = 107.666668
$ws.Range("H40").Value = 3802.4285
$ws.Range("I40").Value = 3035.8572
$ws.Range("K40").Value = 3035.8572
$ws.Range("M40").Value = -2860.8572
$ws.Range("H55").Value = 1359.4117
$ws.Range("I55").Value = 214.875
$ws.Range("J55").Value = 2376.7778
$ws.Range("K55").Value = 214.875
$ws.Range("L55").Value = 2376.7778
$ws.Range("M55").Value = -0.875
$ws.Range("N55").Value = -2804.7778
$ws.Range("H64").Value = 4923.077
$ws.Range("H67").Value = 4923.077
$ws.Range("H74").Value = 3617.875
$ws.Range("I74").Value = 3941.8572
$ws.Range("K74").Value = 3941.8572
$ws.Range("M74").Value = -3005.8572
$ws.Range("H77").Value = 3617.875
$ws.Range("I77").Value = 3941.8572
$ws.Range("K77").Value = 19709.286
$ws.Range("M77").Value = -15029.286
$ws.Range("H80").Value = 729
$ws.Range("I80").Value = 422.5
$ws.Range("J80").Value = 933.3333
$ws.Range("K80").Value = 1267.5
$ws.Range("L80").Value = 2799.9999
$ws.Range("M80").Value = -269.5
$ws.Range("N80").Value = -4795.9999
$ws.Range("H83").Value = 729
$ws.Range("I83").Value = 422.5
$ws.Range("J83").Value = 933.3333
$ws.Range("K83").Value = 3802.5
$ws.Range("L83").Value = 8399.9997
$ws.Range("M83").Value = 1189.5
$ws.Range("N83").Value = -18383.9997
$ws.Range("I113").Value = 50001170
$ws.Range("J113").Value = 133336100
$ws.Range("K113").Value = 50001170
$ws.Range("L113").Value = 133336100
$ws.Range("M113").Value = -49997916
$ws.Range("N113").Value = -133342608
$ws.Range("H132").Value = 1416.2632
$ws.Range("I132").Value = 1416.2632
$ws.Range("K132").Value = 4248.7896
$ws.Range("M132").Value = -1718.7896
$ws.Range("H137").Value = 4986.1
$ws.Range("I137").Value = 3125
$ws.Range("J137").Value = 7260.778
$ws.Range("K137").Value = 9375
$ws.Range("L137").Value = 21782.334
$ws.Range("M137").Value = -6825
$ws.Range("N137").Value = -26882.334
$ws.Range("H138").Value = 1012464.5
$ws.Range("I138").Value = 827.9091
$ws.Range("J138").Value = 1518282.8
$ws.Range("K138").Value = 2483.7273
$ws.Range("L138").Value = 4554848.4
$ws.Range("M138").Value = 2656.2727
$ws.Range("N138").Value = -4565128.4

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11905753
$ws.Range("J32").Value = 4600
$ws.Range("L32").Value = 4600
$ws.Range("N32").Value = -5174
$ws.Range("H61").Value = 45551070
$ws.Range("I61").Value = 83334584
$ws.Range("J61").Value = 210861.6
$ws.Range("K61").Value = 83334584
$ws.Range("L61").Value = 210861.6
$ws.Range("M61").Value = -83334372
$ws.Range("N61").Value = -211285.6
$ws.Range("H88").Value = 1564.2354
$ws.Range("I88").Value = 1428.7778
$ws.Range("J88").Value = 1716.625
$ws.Range("K88").Value = 1428.7778
$ws.Range("L88").Value = 1716.625
$ws.Range("M88").Value = -1022.7778
$ws.Range("N88").Value = -2528.625
$ws.Range("H91").Value = 1564.2354
$ws.Range("I91").Value = 1428.7778
$ws.Range("J91").Value = 1716.625
$ws.Range("K91").Value = 1428.7778
$ws.Range("L91").Value = 1716.625
$ws.Range("M91").Value = -24.77780000000007
$ws.Range("N91").Value = -4524.625
$ws.Range("H92").Value = 84989
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").Value = $null
$ws.Range("H122").Value = 4248.875
$ws.Range("I122").Value = 1996
$ws.Range("J122").Value = 4999.8335
$ws.Range("K122").Value = 5988
$ws.Range("L122").Value = 14999.5005
$ws.Range("M122").Value = -3538
$ws.Range("N122").Value = -19899.5005
$ws.Range("H132").Value = 5745.943
$ws.Range("I132").Value = 3378.577
$ws.Range("J132").Value = 12585
$ws.Range("K132").Value = 10135.731
$ws.Range("L132").Value = 37755
$ws.Range("M132").Value = -7605.731
$ws.Range("N132").Value = -42815
$ws.Range("H136").Value = 45551070
$ws.Range("I136").Value = 83334584
$ws.Range("J136").Value = 210861.6
$ws.Range("K136").Value = 250003752
$ws.Range("L136").Value = 632584.8
$ws.Range("M136").Value = -250001202
$ws.Range("N136").Value = -637684.8

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2256.6924
$ws.Range("I86").Value = 2334.1
$ws.Range("K86").Value = 2334.1
$ws.Range("M86").Value = -1211.1
$ws.Range("H89").Value = 2256.6924
$ws.Range("I89").Value = 2334.1
$ws.Range("K89").Value = 11670.5
$ws.Range("M89").Value = -6054.5
$ws.Range("H102").Value = 88181.5
$ws.Range("I102").Value = 73518.664
$ws.Range("K102").Value = 73518.664
$ws.Range("M102").Value = -70273.664
$ws.Range("H124").Value = 81663.336
$ws.Range("J124").Value = 81663.336
$ws.Range("L124").Value = 81663.336
$ws.Range("N124").Value = -91483.336

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 2623.0588
$ws.Range("I7").Value = 416.6
$ws.Range("K7").Value = 416.6
$ws.Range("M7").Value = -303.6
$ws.Range("H22").Value = 600
$ws.Range("I22").Value = 600
$ws.Range("K22").Value = 600
$ws.Range("M22").Value = -250

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 17466.7
$ws.Range("J44").Value = 126.333336
$ws.Range("L44").Value = 379.000008
$ws.Range("N44").Value = -1175.000008
$ws.Range("H113").Value = 1350.5555
$ws.Range("I113").Value = 500
$ws.Range("J113").Value = 1456.875
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 4370.625
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -8710.625
$ws.Range("H121").Value = 4098.6665
$ws.Range("I121").Value = 1029.5
$ws.Range("J121").Value = 5633.25
$ws.Range("K121").Value = 3088.5
$ws.Range("L121").Value = 16899.75
$ws.Range("M121").Value = -1778.5
$ws.Range("N121").Value = -19519.75

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3609.3125
$ws.Range("I113").Value = 2541.5
$ws.Range("J113").Value = 4250
$ws.Range("K113").Value = 2541.5
$ws.Range("L113").Value = 4250
$ws.Range("M113").Value = -371.5
$ws.Range("N113").Value = -8590

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 54970.4
$ws.Range("J7").Value = 96501.91
$ws.Range("L7").Value = 96501.91
$ws.Range("N7").Value = -96725.91
$ws.Range("H16").Value = 1183.1177
$ws.Range("I16").Value = 1317.4286
$ws.Range("K16").Value = 1317.4286
$ws.Range("M16").Value = -1147.4286
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = $null
$ws.Range("N35").Value = $null
$ws.Range("H98").Value = 74565
$ws.Range("J98").Value = 74565
$ws.Range("L98").Value = 74565
$ws.Range("N98").Value = -80555
$ws.Range("H122").Value = 5060.1177
$ws.Range("I122").Value = 4439.4585
$ws.Range("K122").Value = 13318.3755
$ws.Range("M122").Value = -10868.3755
$ws.Range("H126").Value = 54970.4
$ws.Range("J126").Value = 96501.91
$ws.Range("L126").Value = 289505.73
$ws.Range("N126").Value = -294445.73

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5515.16
$ws.Range("J122").Value = 9836.5
$ws.Range("L122").Value = 29509.5
$ws.Range("N122").Value = -34409.5
$ws.Range("H126").Value = 14068.25
$ws.Range("I126").Value = 11868.223
$ws.Range("K126").Value = 35604.669
$ws.Range("M126").Value = -33134.669

Write-Host "Applied 197 cell updates across 8 sheets"
